$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain text so Excel does not coerce
# numeric-looking strings (e.g. "569.93") into real numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '64.380.13'
$ws.Range("E2").Value = '  -2.90%  '

$ws.Range("D3").Value = '3.170.91'
$ws.Range("E3").Value = '  -4.56%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '569.93'
$ws.Range("E5").Value = '  -2.72%  '

$ws.Range("D6").Value = '168.39'
$ws.Range("E6").Value = '  -8.10%  '

$ws.Range("D7").Value = '0.607'
$ws.Range("E7").Value = '  -6.02%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").Value = '3.173.30'
$ws.Range("E9").Value = '  -4.44%  '

$ws.Range("E10").Value = '  -4.78%  '

$ws.Range("D11").Value = '6.76'
$ws.Range("E11").Value = '  -0.18%  '

$ws.Range("D12").Value = '0.385'
$ws.Range("E12").Value = '  -4.21%  '

$ws.Range("D13").Value = '3.725.64'
$ws.Range("E13").Value = '  -4.45%  '

$ws.Range("D14").Value = '0.128'
$ws.Range("E14").Value = '  -2.52%  '

$ws.Range("D15").Value = '64.430.41'
$ws.Range("E15").Value = '  -2.88%  '

$ws.Range("D16").Value = '25.36'
$ws.Range("E16").Value = '  -3.81%  '

$ws.Range("E17").Value = '  -3.06%  '

$ws.Range("D18").Value = '3.166.83'
$ws.Range("E18").Value = '  -5.63%  '

$ws.Range("D19").Value = '418.00'
$ws.Range("E19").Value = '  -2.83%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '5.37'
$ws.Range("E20").Value = '  -2.98%  '

$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '12.83'
$ws.Range("E21").Value = '  -3.58%  '

$ws.Range("D22").Value = '7.05'
$ws.Range("E22").Value = '  -5.02%  '

$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("D24").Value = '69.84'
$ws.Range("E24").Value = '  -3.18%  '

$ws.Range("E25").Value = '  +0.53%  '

$ws.Range("E26").Value = '  -5.92%  '

$ws.Range("D27").Value = '0.0000106'
$ws.Range("E27").Value = '  -7.46%  '

$ws.Range("D28").Value = '8.80'
$ws.Range("E28").Value = '  -2.16%  '

$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -1.08%  '

$ws.Range("E30").Value = '  -5.90%  '

$ws.Range("D31").Value = '21.68'
$ws.Range("E31").Value = '  -3.11%  '

$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("E33").Value = '  -3.51%  '

$ws.Range("D34").Value = '6.31'
$ws.Range("E34").Value = '  -4.68%  '

$ws.Range("E35").Value = '  -4.73%  '

$ws.Range("D36").Value = '157.24'
$ws.Range("E36").Value = '  -1.65%  '

$ws.Range("E37").Value = '  -6.56%  '

$ws.Range("D38").Value = '2.731.98'
$ws.Range("E38").Value = '  -5.80%  '

$ws.Range("D39").Value = '1.69'
$ws.Range("E39").Value = '  -6.93%  '

$ws.Range("D40").Value = '24.25'
$ws.Range("E40").Value = '  -9.09%  '

$ws.Range("E41").Value = '  -3.97%  '

$ws.Range("E42").Value = '  -2.60%  '

$ws.Range("E43").Value = '  -7.78%  '

$ws.Range("D44").Value = '0.0618'
$ws.Range("E44").Value = '  -7.45%  '

$ws.Range("D45").Value = '5.58'
$ws.Range("E45").Value = '  -7.19%  '

$ws.Range("E46").Value = '  -4.07%  '

$ws.Range("E47").Value = '  -7.91%  '

$ws.Range("D48").Value = '292.44'
$ws.Range("E48").Value = '  -7.71%  '

$ws.Range("E49").Value = '  +0.01%  '

$ws.Range("E50").Value = '  -13.76%  '

$ws.Range("D51").Value = '0.0985'
$ws.Range("E51").Value = '  -6.17%  '
